$wb = $excel.ActiveWorkbook

# --- Add the new "Rekening" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Rekening"

# --- Header row ---
$ws.Range("A1").Value = "No Rekening"
$ws.Range("B1").Value = "Nama Pemilik"
$ws.Range("C1").Value = "Saldo"
$ws.Range("A1:C1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1:C1").VerticalAlignment = -4108     # xlCenter

# --- "No Rekening" column (fill before "Nama Pemilik" so shared strings line up) ---
$ws.Range("A2").Value = "05-0311-224"
$ws.Range("A3").Value = "05-0645-714"
$ws.Range("A4").Value = "04-7653-991"
$ws.Range("A5").Value = "01-4563-202"

# --- "Nama Pemilik" column ---
$ws.Range("B2").Value = "Afif"
$ws.Range("B3").Value = "Ilham"
$ws.Range("B4").Value = "Naufal"
$ws.Range("B5").Value = "Rizki"

# --- "Saldo" column ---
$ws.Range("C2").Value = 10000000
$ws.Range("C3").Value = 10000000
$ws.Range("C4").Value = 10000000
$ws.Range("C5").Value = 10000000
$ws.Range("C2:C5").NumberFormat = "_-* #,##0_-;\-* #,##0_-;_-* ""-""_-;_-@_-"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 17.140625
$ws.Columns.Item(2).ColumnWidth = 14.42578125
$ws.Columns.Item(3).ColumnWidth = 12.5703125

# --- Selection on the new sheet (matches the authored file's cursor position) ---
$ws.Range("E5").Select() | Out-Null
